# "Added frozen headers for printing"
#
# - Freeze the top two header rows so they stay visible while scrolling,
#   with the window scrolled down toward row 105 and the selection left
#   on A3 (bottom/unfrozen pane).
# - Repeat rows 1:2 on every printed page (Print Titles).
# - Re-touch the six left-hand merged header cells (unmerge + re-merge)
#   so they re-serialize after the other six already present, matching
#   the reordering of <mergeCells> in the target file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# --- mergeCells reorder -------------------------------------------------
$headerMerges = @("M1:N1", "C1:D1", "E1:F1", "G1:H1", "I1:J1", "K1:L1")
foreach ($ref in $headerMerges) {
    $ws.Range($ref).UnMerge() | Out-Null
}
foreach ($ref in $headerMerges) {
    $ws.Range($ref).Merge() | Out-Null
}

# --- Print Titles: repeat rows 1-2 on every printed page -----------------
$ws.PageSetup.PrintTitleRows = '$1:$2'

# --- Freeze panes: lock rows 1-2, scroll down, select A3 ------------------
$ws.Range("A3").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# Scroll the window so row 105 is near the top of the unfrozen area.
$excel.ActiveWindow.ScrollRow = 105
$excel.ActiveWindow.ScrollColumn = 1

# Final selection stays on A3, matching the saved view state.
$ws.Range("A3").Select() | Out-Null
